$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$rows = @(491,492,493,494,495,496,497,498,499,500,502,565,579,581)

foreach ($r in $rows) {
    $ws.Range("I$r").Value = "`$1,000,001 - `$5,000,000"
}
